# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from 2023-09-16 (45185) to 2023-10-05 (45204).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
